$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "depara" rows appended to the lookup table (tidy armas rio).
# Columns: A = tipo, B = tipo_formatado, C = flag_arma
# flag is either "TRUE" or "FALSE" -- these must stay literal text (shared
# string), matching the existing table, not COM's auto-boolean coercion.
# We therefore seed each new row by copying an existing row with the same
# flag (which carries the correct style + shared-string text type) and then
# overwrite the A/B text afterwards.
$newRows = @(
    @(24, "revólver", "revolver", "TRUE"),
    @(25, "garrucha", "garrucha", "TRUE"),
    @(26, "pistola", "pistola", "TRUE"),
    @(27, "espingarda", "espingarda", "TRUE"),
    @(28, "fuzil", "fuzil", "TRUE"),
    @(29, "carabina", "carabina", "TRUE"),
    @(30, "submetralhadora", "submetralhadora", "TRUE"),
    @(31, "arma de fabricação caseira", "artesanal", "TRUE"),
    @(32, "outros", $null, "FALSE"),
    @(33, "metralhadora", "metralhadora", "TRUE"),
    @(34, "garruchão", "espingarda", "TRUE")
)

# Row 3 is a full A:C "TRUE" template, row 2 is a full A:C "FALSE" template
# (row 2 has no value in column B, matching the new "outros" row's shape).
$trueTemplate = "3"
$falseTemplate = "2"

foreach ($item in $newRows) {
    $r = $item[0]
    $a = $item[1]
    $b = $item[2]
    $c = $item[3]

    if ($c -eq "TRUE") {
        $tmpl = $trueTemplate
    } else {
        $tmpl = $falseTemplate
    }

    if ($null -eq $b) {
        # Skip column B entirely -- no cell should be written there at all.
        $ws.Range("A$tmpl").Copy($ws.Range("A$r"))
        $ws.Range("C$tmpl").Copy($ws.Range("C$r"))
    } else {
        $ws.Range("A${tmpl}:C${tmpl}").Copy($ws.Range("A${r}:C${r}"))
    }

    $ws.Range("A$r").Value = $a
    if ($null -ne $b) {
        $ws.Range("B$r").Value = $b
    }
}

# Column A now needs to fit the longer labels (matches target bestFit width
# of 21.5 once stored -- the engine re-derives the stored width from the
# visual character width via the default font's max-digit-width, so we feed
# it the pre-image that round-trips to exactly 21.5).
$ws.Columns.Item(1).ColumnWidth = 20.666666666666668

# Adjust the visible view to match the target workbook (scrolled down, new selection).
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("E32").Select()
